# Updates Leve profit-tracking numbers (currentAveragePrice* / LevePrice* /
# LeveProfit*) across all 8 class sheets, per the scheduled-runner refresh.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2: H2, I2, J2, K2, L2, M2, N2
$ws.Range("H2").Value = 368.81818
$ws.Range("I2").Value = 276
$ws.Range("J2").Value = 446.16666
$ws.Range("K2").Value = 276
$ws.Range("L2").Value = 446.16666
$ws.Range("M2").Value = -163
$ws.Range("N2").Value = -672.16666
# Row 33: H33, I33, K33, M33
$ws.Range("H33").Value = 160.27777
$ws.Range("I33").Value = 167.35294
$ws.Range("K33").Value = 167.35294
$ws.Range("M33").Value = 61.64706000000001
# Row 69: H69, J69, L69, N69
$ws.Range("H69").Value = 1561.1111
$ws.Range("J69").Value = 1532.6923
$ws.Range("L69").Value = 4598.0769
$ws.Range("N69").Value = -6346.0769
# Row 72: H72, J72, L72, N72
$ws.Range("H72").Value = 1561.1111
$ws.Range("J72").Value = 1532.6923
$ws.Range("L72").Value = 13794.2307
$ws.Range("N72").Value = -22530.2307
# Row 98: H98, I98, J98, K98, L98, M98, N98
$ws.Range("H98").Value = 843.5161000000001
$ws.Range("I98").Value = 689.95654
$ws.Range("J98").Value = 1285
$ws.Range("K98").Value = 689.95654
$ws.Range("L98").Value = 1285
$ws.Range("M98").Value = 808.04346
$ws.Range("N98").Value = -4281
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 843.5161000000001
$ws.Range("I122").Value = 689.95654
$ws.Range("J122").Value = 1285
$ws.Range("K122").Value = 2069.86962
$ws.Range("L122").Value = 3855
$ws.Range("M122").Value = 380.1303800000001
$ws.Range("N122").Value = -8755
# Row 129: H129, I129, J129, K129, L129, M129, N129
$ws.Range("H129").Value = 164845.19
$ws.Range("I129").Value = 231.22223
$ws.Range("J129").Value = 193336.06
$ws.Range("K129").Value = 693.66669
$ws.Range("L129").Value = 580008.1799999999
$ws.Range("M129").Value = 4306.33331
$ws.Range("N129").Value = -590008.1799999999
# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 50780
$ws.Range("J140").Value = 50780
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2: H2, I2, J2, K2, L2, M2, N2
$ws.Range("H2").Value = 664.64514
$ws.Range("I2").Value = 647.9655
$ws.Range("J2").Value = 906.5
$ws.Range("K2").Value = 647.9655
$ws.Range("L2").Value = 906.5
$ws.Range("M2").Value = -534.9655
$ws.Range("N2").Value = -1132.5
# Row 74: H74, I74, K74, M74
$ws.Range("H74").Value = 35716252
$ws.Range("I74").Value = 71429280
$ws.Range("K74").Value = 71429280
$ws.Range("M74").Value = -71428406
# Row 77: H77, I77, K77, M77
$ws.Range("H77").Value = 35716252
$ws.Range("I77").Value = 71429280
$ws.Range("K77").Value = 357146400
$ws.Range("M77").Value = -357142032
# Row 97: H97, I97, K97, M97
$ws.Range("H97").Value = 62501212
$ws.Range("I97").Value = 1046.6923
$ws.Range("K97").Value = 1046.6923
$ws.Range("M97").Value = -550.6922999999999
# Row 116: H116, I116, J116, K116, L116, M116, N116
$ws.Range("H116").Value = 664.64514
$ws.Range("I116").Value = 647.9655
$ws.Range("J116").Value = 906.5
$ws.Range("K116").Value = 647.9655
$ws.Range("L116").Value = 906.5
$ws.Range("M116").Value = 1646.0345
$ws.Range("N116").Value = -5494.5
# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3: H3, I3, J3, K3, L3, M3, N3
$ws.Range("H3").Value = 664.64514
$ws.Range("I3").Value = 647.9655
$ws.Range("J3").Value = 906.5
$ws.Range("K3").Value = 647.9655
$ws.Range("L3").Value = 906.5
$ws.Range("M3").Value = -533.9655
$ws.Range("N3").Value = -1134.5
# Row 11: H11, I11, K11, M11
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 1000
$ws.Range("M11").Value = -860
# Row 95: H95, J95, L95, N95
$ws.Range("H95").Value = 19999.666
$ws.Range("J95").Value = 19999.666
$ws.Range("L95").Value = 19999.666
$ws.Range("N95").Value = -25491.666
# Row 99: H99, I99, K99, M99
$ws.Range("H99").Value = 1970.3334
$ws.Range("I99").Value = 1950
$ws.Range("K99").Value = 1950
$ws.Range("M99").Value = -452
# Row 105: H105, I105, J105, K105, L105, M105, N105
$ws.Range("H105").Value = 1516831.5
$ws.Range("I105").Value = 1402.8
$ws.Range("J105").Value = 2175713.5
$ws.Range("K105").Value = 1402.8
$ws.Range("L105").Value = 2175713.5
$ws.Range("M105").Value = 344.2
$ws.Range("N105").Value = -2179207.5
# Row 107: H107, I107, K107, M107
$ws.Range("H107").Value = 1414.421
$ws.Range("I107").Value = 1306.3334
$ws.Range("K107").Value = 1306.3334
$ws.Range("M107").Value = 613.6666
# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 3970.3333
$ws.Range("I134").Value = 4018.5312
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 12055.5936
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -9520.5936
$ws.Range("N134").Value = -16320

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16: H16, I16, J16, K16, L16, M16, N16
$ws.Range("H16").Value = 1196.6666
$ws.Range("I16").Value = 1196.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1196.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -909.6666
$ws.Range("N16").ClearContents()
# Row 105: H105, I105, K105, M105
$ws.Range("H105").Value = 1036.1052
$ws.Range("I105").Value = 927.94116
$ws.Range("K105").Value = 927.94116
$ws.Range("M105").Value = 819.05884
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 1263.56
$ws.Range("I107").Value = 475.35715
$ws.Range("J107").Value = 2266.7273
$ws.Range("K107").Value = 475.35715
$ws.Range("L107").Value = 2266.7273
$ws.Range("M107").Value = 1444.64285
$ws.Range("N107").Value = -6106.7273
# Row 109: H109, J109, L109, N109
$ws.Range("H109").Value = 198020380
$ws.Range("J109").Value = 198020380
$ws.Range("L109").Value = 198020380
$ws.Range("N109").Value = -198022460
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 1196.6666
$ws.Range("I113").Value = 1196.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1196.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 973.3334
$ws.Range("N113").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 122: H122, J122, L122, N122
$ws.Range("H122").Value = 612.61536
$ws.Range("J122").Value = 681.0526
$ws.Range("L122").Value = 6129.4734
$ws.Range("N122").Value = -11029.4734
# Row 131: H131, J131, L131, N131
$ws.Range("H131").Value = 695.92
$ws.Range("J131").Value = 724.0879
$ws.Range("L131").Value = 2172.2637
$ws.Range("N131").Value = -12252.2637
# Row 132: H132, J132, L132, N132
$ws.Range("H132").Value = 585
$ws.Range("J132").Value = 568.3333
$ws.Range("L132").Value = 5114.9997
$ws.Range("N132").Value = -10174.9997

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 15: H15, J15, L15, N15
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20576
# Row 81: H81, J81, L81, N81
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996
# Row 84: H84, J84, L84, N84
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 7937.3125
$ws.Range("I113").Value = 9446.416999999999
$ws.Range("J113").Value = 3410
$ws.Range("K113").Value = 9446.416999999999
$ws.Range("L113").Value = 3410
$ws.Range("M113").Value = -7276.416999999999
$ws.Range("N113").Value = -7750
# Row 136: H136, J136, L136, N136
$ws.Range("H136").Value = 9919.923000000001
$ws.Range("J136").Value = 9919.923000000001
$ws.Range("L136").Value = 29759.769
$ws.Range("N136").Value = -34859.769
# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 55000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360
# Row 141: H141, J141, L141, N141
$ws.Range("H141").Value = 44452.867
$ws.Range("J141").Value = 44452.867
$ws.Range("L141").Value = 44452.867
$ws.Range("N141").Value = -54812.867

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 43: H43, J43, L43, N43
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30386
# Row 61: H61, I61, K61, M61
$ws.Range("H61").Value = 4870.7856
$ws.Range("I61").Value = 1721.3334
$ws.Range("K61").Value = 1721.3334
$ws.Range("M61").Value = -1519.3334
# Row 93: H93, I93, K93, M93
$ws.Range("H93").Value = 1104.6154
$ws.Range("I93").Value = 1030
$ws.Range("K93").Value = 1030
$ws.Range("M93").Value = 218
# Row 113: H113, I113, K113, M113
$ws.Range("H113").Value = 4870.7856
$ws.Range("I113").Value = 1721.3334
$ws.Range("K113").Value = 1721.3334
$ws.Range("M113").Value = 448.6666

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 27: H27, J27, L27, N27
$ws.Range("H27").Value = 35350
$ws.Range("J27").Value = 35350
$ws.Range("L27").Value = 35350
$ws.Range("N27").Value = -35488
# Row 115: H115
$ws.Range("H115").Value = 28249.166
# Row 140: H140, J140, L140, N140
$ws.Range("H140").Value = 45319.8
$ws.Range("J140").Value = 45319.8
$ws.Range("L140").Value = 45319.8
$ws.Range("N140").Value = -55679.8
# Row 141: H141, J141, L141, N141
$ws.Range("H141").Value = 72197.8
$ws.Range("J141").Value = 72197.8
$ws.Range("L141").Value = 72197.8
$ws.Range("N141").Value = -82557.8

Write-Output "Applied all Typhon_Profits cell updates."
